$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText "D2" '27.349.57'
Set-CellText "E2" '  -1.79%  '

Set-CellText "D3" '1.727.68'
Set-CellText "E3" '  -2.12%  '

Set-CellText "E4" '  +0.50%  '

Set-CellText "D5" '321.57'
Set-CellText "E5" '  -0.29%  '

Set-CellText "D6" '1.005'
Set-CellText "E6" '  +0.48%  '

Set-CellText "D7" '0.4526'
Set-CellText "E7" '  +6.12%  '

Set-CellText "D8" '0.3516'
Set-CellText "E8" '  -3.31%  '

Set-CellText "D9" '41.51'
Set-CellText "E9" '  -2.85%  '

Set-CellText "D10" '0.07321'
Set-CellText "E10" '  -3.01%  '

Set-CellText "D11" '1.071'
Set-CellText "E11" '  -2.12%  '

Set-CellText "D12" '1.005'
Set-CellText "E12" '  +0.51%  '

Set-CellText "D13" '20.28'
Set-CellText "E13" '  -2.41%  '

Set-CellText "D14" '5.884'
Set-CellText "E14" '  -3.25%  '

Set-CellText "D15" '7.029'
Set-CellText "E15" '  -3.57%  '

Set-CellText "D16" '1.734.54'
Set-CellText "E16" '  -1.05%  '

Set-CellText "D18" '0.00001046'
Set-CellText "E18" '  -1.86%  '

Set-CellText "D19" '0.06326'
Set-CellText "E19" '  -1.01%  '

Set-CellText "D20" '1.004'
Set-CellText "E20" '  +0.49%  '

Set-CellText "D21" '16.53'
Set-CellText "E21" '  -3.24%  '

Set-CellText "D22" '5.732'
Set-CellText "E22" '  -3.20%  '

Set-CellText "D23" '27.408.02'
Set-CellText "E23" '  -1.64%  '

Set-CellText "D24" '11.00'
Set-CellText "E24" '  -2.43%  '

Set-CellText "D25" '2.100'
Set-CellText "E25" '  -0.98%  '

Set-CellText "D26" '162.23'
Set-CellText "E26" '  +1.71%  '

Set-CellText "D27" '19.75'
Set-CellText "E27" '  -2.65%  '

Set-CellText "D28" '1.931.21'
Set-CellText "E28" '  -1.57%  '

Set-CellText "B29" 'LidoDAOToken'
Set-CellText "C29" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText "D29" '2.039'
Set-CellText "E29" '  -5.34%  '

Set-CellText "B30" 'BitcoinCash'
Set-CellText "C30" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText "D30" '124.17'
Set-CellText "E30" '  -0.85%  '

Set-CellText "D31" '1.040'
Set-CellText "E31" '  -7.55%  '

Set-CellText "D32" '0.09096'
Set-CellText "E32" '  +2.39%  '

Set-CellText "D33" '3.657'
Set-CellText "E33" '  -0.75%  '

Set-CellText "D34" '5.321'
Set-CellText "E34" '  -4.63%  '

Set-CellText "B35" 'Aptos'
Set-CellText "C35" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText "D35" '11.58'
Set-CellText "E35" '  -5.76%  '

Set-CellText "B36" 'VeChain'
Set-CellText "C36" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText "D36" '0.02253'
Set-CellText "E36" '  -1.89%  '

Set-CellText "D37" '0.05940'
Set-CellText "E37" '  -1.65%  '

Set-CellText "D38" '0.2040'
Set-CellText "E38" '  -3.36%  '

Set-CellText "D39" '0.6190'
Set-CellText "E39" '  -2.61%  '

Set-CellText "D40" '4.828'
Set-CellText "E40" '  -3.29%  '

Set-CellText "D41" '1.180'
Set-CellText "E41" '  +0.26%  '

Set-CellText "E42" '  -2.05%  '

Set-CellText "D43" '7.673'
Set-CellText "E43" '  -2.86%  '

Set-CellText "D44" '12.96'
Set-CellText "E44" '  -3.15%  '

Set-CellText "D45" '3.686'
Set-CellText "E45" '  -0.13%  '

Set-CellText "D46" '0.5772'
Set-CellText "E46" '  -2.00%  '

Set-CellText "D47" '121.57'
Set-CellText "E47" '  -1.17%  '

Set-CellText "D48" '1.912'
Set-CellText "E48" '  -4.10%  '

Set-CellText "D49" '0.06817'
Set-CellText "E49" '  -0.38%  '

Set-CellText "D50" '1.107'
Set-CellText "E50" '  -6.73%  '

Set-CellText "D51" '70.66'
Set-CellText "E51" '  -4.47%  '
